$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert a new "Test 3" column before the existing "Total" column ---
# Move the "Total" header (with its style) from G1 into the new H1 cell, then
# relabel G1 as "Test 3".
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Total"
$ws.Range("G1").Value = "Test 3"

# --- Data rows: refresh Test 1 (E) / Test 2 (F) rounding, add Test 3 (G) scores,
#     and populate the new Total (H) column. ---
$rows = @(
    @{ Row=2; E=20.2; F=18.96; G=21.25; H=60.42 },
    @{ Row=3; E=13.47; F=17.44; G=18.75; H=49.66 },
    @{ Row=4; E=18.37; F=15.56; G=43.75; H=77.68000000000001 },
    @{ Row=5; E=18.37; F=16.36; G=18.75; H=53.48 },
    @{ Row=6; E=15.31; F=0; G=17.5; H=32.81 },
    @{ Row=7; E=12.25; F=14.09; G=32.5; H=58.83 },
    @{ Row=8; E=16.53; F=16.72; G=31.25; H=64.5 },
    @{ Row=9; E=18.98; F=16.4; G=43.75; H=79.13 },
    @{ Row=10; E=15.31; F=13.68; G=47.5; H=76.48999999999999 },
    @{ Row=11; E=17.75; F=9.04; G=41.25; H=68.04000000000001 },
    @{ Row=12; E=19.59; F=16.88; G=32.5; H=68.97 },
    @{ Row=13; E=26.33; F=17; G=50; H=93.33 },
    @{ Row=14; E=14.08; F=16.68; G=41.25; H=72.01000000000001 },
    @{ Row=15; E=22.04; F=14.13; G=25; H=61.17 },
    @{ Row=16; E=0; F=15.72; G=35; H=50.72 },
    @{ Row=17; E=20.82; F=14.68; G=20; H=55.5 },
    @{ Row=18; E=22.65; F=17.65; G=50; H=90.3 },
    @{ Row=19; E=19.59; F=14.68; G=0; H=34.27 },
    @{ Row=20; E=18.98; F=15; G=33.75; H=67.73 },
    @{ Row=21; E=15.31; F=16.02; G=43.75; H=75.06999999999999 },
    @{ Row=22; E=12.86; F=0; G=46.25; H=59.11 },
    @{ Row=23; E=23.88; F=17.65; G=32.5; H=74.03 },
    @{ Row=24; E=15.92; F=14.36; G=41.25; H=71.53 },
    @{ Row=25; E=22.04; F=16.88; G=22.5; H=61.42 },
    @{ Row=26; E=19.59; F=14.36; G=30; H=63.95 },
    @{ Row=27; E=12.86; F=18.08; G=47.5; H=78.44 },
    @{ Row=28; E=23.26; F=14.36; G=31.25; H=68.88 },
    @{ Row=29; E=20.82; F=15.92; G=16.25; H=52.99 },
    @{ Row=30; E=0; F=0; G=0; H=0 },
    @{ Row=31; E=18.37; F=11; G=28.75; H=58.12 },
    @{ Row=32; E=18.37; F=16.8; G=36.25; H=71.42 },
    @{ Row=33; E=22.65; F=17.76; G=16.25; H=56.66 },
    @{ Row=34; E=16.53; F=10.28; G=0; H=26.81 },
    @{ Row=35; E=17.75; F=15.04; G=25; H=57.79 },
    @{ Row=36; E=20.82; F=13.82; G=0; H=34.63 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
